$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "#10-파이썬(Python) 내장함수-map, zip, filter, enumerate"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-10"

$ws.Range("D5").Value = "오일러-코시 방정식"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/06/07/Euler_Cauchy.html"

$ws.Range("D12").Value = "“혼자 공부하는 머신러닝+딥러닝”, “머신 러닝 교과서 3판”, “딥러닝 일러스트레이티드” 주피터 노트북 텐서플로 2.5 테스트 완료"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/06/03/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d-%eb%a8%b8%ec%8b%a0-%eb%9f%ac%eb%8b%9d-%ea%b5%90%ea%b3%bc%ec%84%9c-3/"
